$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "26.999.76"
$ws.Cells.Item(2, 5).Value = "  -1.92%  "
$ws.Cells.Item(3, 4).Value = "1.820.87"
$ws.Cells.Item(3, 5).Value = "  -0.90%  "
$ws.Cells.Item(4, 4).Value = "'1.001"
$ws.Cells.Item(4, 5).Value = "  -0.44%  "
$ws.Cells.Item(5, 4).Value = "'310.92"
$ws.Cells.Item(5, 5).Value = "  -1.64%  "
$ws.Cells.Item(6, 4).Value = "'1.001"
$ws.Cells.Item(6, 5).Value = "  -0.38%  "
$ws.Cells.Item(7, 4).Value = "'0.4240"
$ws.Cells.Item(7, 5).Value = "  -1.48%  "
$ws.Cells.Item(8, 4).Value = "'0.3659"
$ws.Cells.Item(9, 4).Value = "'0.07208"
$ws.Cells.Item(9, 5).Value = "  -0.98%  "
$ws.Cells.Item(10, 4).Value = "'0.8386"
$ws.Cells.Item(10, 5).Value = "  -3.46%  "
$ws.Cells.Item(11, 4).Value = "'20.55"
$ws.Cells.Item(11, 5).Value = "  -3.30%  "
$ws.Cells.Item(12, 4).Value = "1.828.80"
$ws.Cells.Item(12, 5).Value = "  -0.98%  "
$ws.Cells.Item(13, 4).Value = "'6.628"
$ws.Cells.Item(13, 5).Value = "  -1.18%  "
$ws.Cells.Item(14, 4).Value = "'0.07059"
$ws.Cells.Item(14, 5).Value = "  -0.33%  "
$ws.Cells.Item(15, 4).Value = "'5.263"
$ws.Cells.Item(15, 5).Value = "  -2.00%  "
$ws.Cells.Item(16, 4).Value = "'89.27"
$ws.Cells.Item(16, 5).Value = "  +0.84%  "
$ws.Cells.Item(17, 4).Value = "'1.002"
$ws.Cells.Item(17, 5).Value = "  -0.63%  "
$ws.Cells.Item(18, 4).Value = "'0.000008718"
$ws.Cells.Item(18, 5).Value = "  -2.52%  "
$ws.Cells.Item(19, 4).Value = "'1.002"
$ws.Cells.Item(19, 5).Value = "  -0.34%  "
$ws.Cells.Item(20, 4).Value = "'14.81"
$ws.Cells.Item(20, 5).Value = "  -3.21%  "
$ws.Cells.Item(21, 4).Value = "27.061.03"
$ws.Cells.Item(21, 5).Value = "  -1.74%  "
$ws.Cells.Item(22, 4).Value = "'5.113"
$ws.Cells.Item(22, 5).Value = "  -1.10%  "
$ws.Cells.Item(23, 4).Value = "'10.77"
$ws.Cells.Item(23, 5).Value = "  -1.94%  "
$ws.Cells.Item(24, 4).Value = "2.053.33"
$ws.Cells.Item(24, 5).Value = "  -0.86%  "
$ws.Cells.Item(25, 5).Value = "  -1.63%  "
$ws.Cells.Item(26, 4).Value = "'150.96"
$ws.Cells.Item(26, 5).Value = "  -1.92%  "
$ws.Cells.Item(27, 4).Value = "'2.231"
$ws.Cells.Item(27, 5).Value = "  +3.31%  "
$ws.Cells.Item(28, 4).Value = "'18.19"
$ws.Cells.Item(28, 5).Value = "  -1.41%  "
$ws.Cells.Item(29, 4).Value = "'5.211"
$ws.Cells.Item(29, 5).Value = "  -1.81%  "
$ws.Cells.Item(30, 4).Value = "'116.59"
$ws.Cells.Item(30, 5).Value = "  -0.68%  "
$ws.Cells.Item(31, 4).Value = "'0.08716"
$ws.Cells.Item(31, 5).Value = "  -1.88%  "
$ws.Cells.Item(32, 4).Value = "'1.171"
$ws.Cells.Item(32, 5).Value = "  -3.33%  "
$ws.Cells.Item(33, 4).Value = "'0.7304"
$ws.Cells.Item(33, 5).Value = "  -5.43%  "
$ws.Cells.Item(34, 4).Value = "'2.904"
$ws.Cells.Item(34, 5).Value = "  +0.26%  "
$ws.Cells.Item(35, 4).Value = "'4.398"
$ws.Cells.Item(35, 5).Value = "  -2.31%  "
$ws.Cells.Item(36, 5).Value = "  -0.42%  "
$ws.Cells.Item(37, 4).Value = "'1.089"
$ws.Cells.Item(37, 5).Value = "  -3.20%  "
$ws.Cells.Item(38, 4).Value = "'0.01937"
$ws.Cells.Item(38, 5).Value = "  -1.17%  "
$ws.Cells.Item(39, 4).Value = "'0.05209"
$ws.Cells.Item(39, 5).Value = "  -1.57%  "
$ws.Cells.Item(40, 4).Value = "'7.246"
$ws.Cells.Item(40, 5).Value = "  +1.19%  "
$ws.Cells.Item(41, 4).Value = "'2.862"
$ws.Cells.Item(41, 5).Value = "  -0.48%  "
$ws.Cells.Item(42, 4).Value = "'0.1683"
$ws.Cells.Item(42, 5).Value = "  +0.16%  "
$ws.Cells.Item(43, 4).Value = "'0.5090"
$ws.Cells.Item(43, 5).Value = "  -0.19%  "
$ws.Cells.Item(44, 4).Value = "'8.512"
$ws.Cells.Item(44, 5).Value = "  -2.47%  "
$ws.Cells.Item(45, 4).Value = "'10.45"
$ws.Cells.Item(45, 5).Value = "  -1.55%  "
$ws.Cells.Item(46, 4).Value = "'1.954"
$ws.Cells.Item(46, 5).Value = "  +6.25%  "
$ws.Cells.Item(47, 4).Value = "'0.4715"
$ws.Cells.Item(47, 5).Value = "  -0.28%  "
$ws.Cells.Item(48, 4).Value = "'105.48"
$ws.Cells.Item(48, 5).Value = "  -1.01%  "
$ws.Cells.Item(49, 5).Value = "  -0.43%  "
$ws.Cells.Item(50, 4).Value = "'0.06312"
$ws.Cells.Item(50, 5).Value = "  -1.87%  "
$ws.Cells.Item(51, 5).Value = "  -2.10%  "
